$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last refreshed" timestamp in the title cell
$ws.Range("A1").Value = "Datos actualizados a 5 de Agosto de 2020 a las 11:56"

# Row 4
$ws.Range("E4").Value = 2275548
$ws.Range("G4").Value = 33
$ws.Range("H4").Value = 160323

# Row 19
$ws.Range("B19").Value = 246674
$ws.Range("C19").Value = 2654
$ws.Range("D19").Value = 141750
$ws.Range("E19").Value = 101657
$ws.Range("G19").Value = 33
$ws.Range("H19").Value = 3267

# Row 22
$ws.Range("B22").Value = 213080
$ws.Range("E22").Value = 9141
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = 9239

# Row 26  # -> Indonesia
$ws.Range("A26").Value = "Indonesia"
$ws.Range("B26").Value = 116871
$ws.Range("C26").Value = 1815
$ws.Range("D26").Value = 73889
$ws.Range("E26").Value = 37530
$ws.Range("G26").Value = 64
$ws.Range("H26").Value = 5452

# Row 27  # -> Filipinas
$ws.Range("A27").Value = "Filipinas"
$ws.Range("B27").Value = 115980
$ws.Range("C27").Value = 3462
$ws.Range("D27").Value = 66270
$ws.Range("E27").Value = 47587
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = 2123

# Row 49
$ws.Range("B49").Value = 48789
$ws.Range("C49").Value = 640
$ws.Range("E49").Value = 11712
$ws.Range("G49").Value = 18
$ws.Range("H49").Value = 1756

# Row 52
$ws.Range("E52").Value = 2644
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 153

# Row 67
$ws.Range("B67").Value = 21566
$ws.Range("C67").Value = 85
$ws.Range("D67").Value = 19464
$ws.Range("E67").Value = 1383

# Row 73
$ws.Range("B73").Value = 18701
$ws.Range("C73").Value = 439
$ws.Range("D73").Value = 9111
$ws.Range("E73").Value = 9092
$ws.Range("G73").Value = 12
$ws.Range("H73").Value = 498

# Row 80
$ws.Range("E80").Value = 6263
$ws.Range("G80").Value = 2
$ws.Range("H80").Value = 88

# Row 88
$ws.Range("B88").Value = 9023
$ws.Range("C88").Value = 21
$ws.Range("D88").Value = 8702
$ws.Range("E88").Value = 196

# Row 93  # -> Finlandia
$ws.Range("A93").Value = "Finlandia"
$ws.Range("B93").Value = 7512
$ws.Range("C93").Value = 29
$ws.Range("D93").Value = 6950
$ws.Range("E93").Value = 231
$ws.Range("H93").Value = 331

# Row 94  # -> Guinea
$ws.Range("A94").Value = "Guinea"
$ws.Range("B94").Value = 7489
$ws.Range("D94").Value = 6591
$ws.Range("E94").Value = 850
$ws.Range("H94").Value = 48

# Row 99
$ws.Range("D99").Value = 3054
$ws.Range("E99").Value = 2520

# Row 101  # -> Libano
$ws.Range("A101").Value = "Libano"
$ws.Range("B101").Value = 5271
$ws.Range("D101").Value = 1837
$ws.Range("E101").Value = 3369
$ws.Range("H101").Value = 65

# Row 102  # -> Republica de Yibuti
$ws.Range("A102").Value = "Republica de Yibuti"
$ws.Range("B102").Value = 5248
$ws.Range("D102").Value = 5044
$ws.Range("E102").Value = 145
$ws.Range("H102").Value = 59

# Row 106
$ws.Range("B106").Value = 4564
$ws.Range("C106").Value = 11
$ws.Range("D106").Value = 3431
$ws.Range("E106").Value = 534
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 599

# Row 112
$ws.Range("B112").Value = 3755
$ws.Range("C112").Value = 85
$ws.Range("D112").Value = 2314
$ws.Range("E112").Value = 1398
$ws.Range("G112").Value = 2
$ws.Range("H112").Value = 43

# Row 119
$ws.Range("D119").Value = 2537
$ws.Range("E119").Value = 286

# Row 126
$ws.Range("B126").Value = 2208
$ws.Range("C126").Value = 18
$ws.Range("D126").Value = 1884
$ws.Range("E126").Value = 200
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 124

# Row 127
$ws.Range("B127").Value = 2147
$ws.Range("C127").Value = 10
$ws.Range("D127").Value = 1650
$ws.Range("E127").Value = 416
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 81
